# NIT-9013779150.xlsx edit: elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Preserve the old "last row" formatting (row 19's current style)
#    by copying it down to the brand-new row 27 BEFORE row 19 itself
#    is overwritten with the "normal" row style.
# ---------------------------------------------------------------------
$ws.Range("B19:J19").Copy($ws.Range("B27:J27"))

# ---------------------------------------------------------------------
# 2) Turn old row 19 into a "normal" data row (same style as row 18),
#    and stamp the same style on the newly inserted rows 20-26.
# ---------------------------------------------------------------------
$ws.Range("B18:J18").Copy($ws.Range("B19:J19"))
$ws.Range("B18:J18").Copy($ws.Range("B20:J20"))
$ws.Range("B18:J18").Copy($ws.Range("B21:J21"))
$ws.Range("B18:J18").Copy($ws.Range("B22:J22"))
$ws.Range("B18:J18").Copy($ws.Range("B23:J23"))
$ws.Range("B18:J18").Copy($ws.Range("B24:J24"))
$ws.Range("B18:J18").Copy($ws.Range("B25:J25"))
$ws.Range("B18:J18").Copy($ws.Range("B26:J26"))

# Clear any stray leftover values in H/I/J for all data rows (they must
# stay blank; only the B:G columns carry data).
$ws.Range("H16:J27").ClearContents()

# ---------------------------------------------------------------------
# 3) Move the signature block (old rows 24-25) down to rows 32-33,
#    preserving their formatting, then blank out the old location.
# ---------------------------------------------------------------------
$ws.Range("B24:J24").Copy($ws.Range("B32:J32"))
$ws.Range("B25:J25").Copy($ws.Range("B33:J33"))
$ws.Range("B24:J25").Clear()

$ws.Range("B24:C24").UnMerge()
$ws.Range("B25:C25").UnMerge()
$ws.Range("H24:J24").UnMerge()
$ws.Range("H25:J25").UnMerge()
$ws.Range("B32:C32").Merge()
$ws.Range("B33:C33").Merge()
$ws.Range("H32:J32").Merge()
$ws.Range("H33:J33").Merge()

# ---------------------------------------------------------------------
# 4) Header / summary cells (text is unchanged, only E11 and F13 carry
#    real numeric changes).
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 407419
$ws.Range("F13").Value = 11

# ---------------------------------------------------------------------
# 5) Worker / debt rows - full refresh of the table contents.
# ---------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73350807"
$ws.Range("D16").Value = "ASMETH LEONAR MARRUGO GONZALEZ"
$ws.Range("E16").Value = "2012"
$ws.Range("F16").Value = 1211
$ws.Range("G16").Value = 908526

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1051444928"
$ws.Range("D17").Value = "LUIS TOMAS CAUSIL PATERNINA"
$ws.Range("E17").Value = "2012"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 908526

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143359583"
$ws.Range("D18").Value = "ROBINSON FELIPE GELVIS PACHECO"
$ws.Range("E18").Value = "2206"
$ws.Range("F18").Value = 27867
$ws.Range("G18").Value = 908526

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143359583"
$ws.Range("D19").Value = "ROBINSON FELIPE GELVIS PACHECO"
$ws.Range("E19").Value = "2205"
$ws.Range("F19").Value = 38000
$ws.Range("G19").Value = 908526

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143359583"
$ws.Range("D20").Value = "ROBINSON FELIPE GELVIS PACHECO"
$ws.Range("E20").Value = "2204"
$ws.Range("F20").Value = 38000
$ws.Range("G20").Value = 908526

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1143359583"
$ws.Range("D21").Value = "ROBINSON FELIPE GELVIS PACHECO"
$ws.Range("E21").Value = "2203"
$ws.Range("F21").Value = 38000
$ws.Range("G21").Value = 908526

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1143359583"
$ws.Range("D22").Value = "ROBINSON FELIPE GELVIS PACHECO"
$ws.Range("E22").Value = "2202"
$ws.Range("F22").Value = 38000
$ws.Range("G22").Value = 908526

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1143359583"
$ws.Range("D23").Value = "ROBINSON FELIPE GELVIS PACHECO"
$ws.Range("E23").Value = "2201"
$ws.Range("F23").Value = 38000
$ws.Range("G23").Value = 908526

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1143359583"
$ws.Range("D24").Value = "ROBINSON FELIPE GELVIS PACHECO"
$ws.Range("E24").Value = "2112"
$ws.Range("F24").Value = 38000
$ws.Range("G24").Value = 908526

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1143359583"
$ws.Range("D25").Value = "ROBINSON FELIPE GELVIS PACHECO"
$ws.Range("E25").Value = "2111"
$ws.Range("F25").Value = 38000
$ws.Range("G25").Value = 908526

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1143359583"
$ws.Range("D26").Value = "ROBINSON FELIPE GELVIS PACHECO"
$ws.Range("E26").Value = "2010"
$ws.Range("F26").Value = 38000
$ws.Range("G26").Value = 908526

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1143359583"
$ws.Range("D27").Value = "ROBINSON FELIPE GELVIS PACHECO"
$ws.Range("E27").Value = "2009"
$ws.Range("F27").Value = 38000
$ws.Range("G27").Value = 908526
